$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64, shifting rows 64..108 down to 65..109.
$ws.Rows.Item(64).Insert()

# Fill the newly-inserted row 64: the "fixed" columns repeat the same
# market/category info found in every other data row, and D/J/K/L/M/O/P
# carry this observation's own values.
$ws.Range("A64").Value = 5
$ws.Range("B64").Value = "Macroferia Regional de Talca"
$ws.Range("C64").Value = "Maule"
$ws.Range("D64").Value = 44879
$ws.Range("E64").Value = 7
$ws.Range("F64").Value = 100112026
$ws.Range("G64").Value = "Haba"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 300
$ws.Range("K64").Value = 9000
$ws.Range("L64").Value = 9000
$ws.Range("M64").Value = 9000
$ws.Range("N64").Value = '$/saco 25 kilos'
$ws.Range("O64").Value = "Región del Maule"
$ws.Range("P64").Value = 360
$ws.Range("Q64").Value = 25
$ws.Range("R64").Value = "Hortaliza"

# Match the date-number-format style used by the other "Fecha" column cells.
$ws.Range("D64").NumberFormat = $ws.Range("D65").NumberFormat
